$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I5").Value = 0.6295521131073972
$ws.Range("J5").Value = 0.4508479156658732
$ws.Range("K5").Value = -0.001005534651262963
$ws.Range("L5").Value = 2.434995522272491
